# database/industries/siman/sesoufi/income/yearly/dollar.xlsx
# "update database and change read_price algorithm"
#
# The yearly income-statement table keeps a rolling 5-year window
# (columns D..H). Refreshing the database drops the oldest reported
# year, slides the remaining four years one column to the left, and
# appends the newly reported year's figures in column H. The period
# headers (row 8) and publish-date stamps (row 9) shift the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "12 ماهه منتهی به ####/12" period headers, shifted left one year ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" publish-date stamps, shifted left one year ---
$ws.Range("D9").Value = "1399-01-27 (10)"
$ws.Range("E9").Value = "1400-02-05 (8)"
$ws.Range("F9").Value = "1401-02-10 (9)"
$ws.Range("G9").Value = "1402-02-12 (10)"
$ws.Range("H9").Value = "1402-02-12 (2)"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 17339
$ws.Range("E11").Value = 21284
$ws.Range("F11").Value = 20107
$ws.Range("G11").Value = 30949
$ws.Range("H11").Value = 34419

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -12674
$ws.Range("E12").Value = -14633
$ws.Range("F12").Value = -13228
$ws.Range("G12").Value = -18884
$ws.Range("H12").Value = -21362

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = 6651
$ws.Range("F13").Value = 6879
$ws.Range("G13").Value = 12065
$ws.Range("H13").Value = 13057

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -1665
$ws.Range("E14").Value = -2355
$ws.Range("F14").Value = -1535
$ws.Range("G14").Value = -1664
$ws.Range("H14").Value = -2108

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (unchanged, still all dashes) ---
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating inc/exp) ---
$ws.Range("D16").Value = -742
$ws.Range("E16").Value = -558
$ws.Range("F16").Value = 363
$ws.Range("G16").Value = -1305
$ws.Range("H16").Value = 323

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = 3739
$ws.Range("F17").Value = 5707
$ws.Range("G17").Value = 9096
$ws.Range("H17").Value = 11271

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = -818
$ws.Range("E18").Value = -521
$ws.Range("F18").Value = -66
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = -21

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Non-operating inc/exp) ---
$ws.Range("D19").Value = -110
$ws.Range("E19").Value = 293
$ws.Range("F19").Value = 1565
$ws.Range("G19").Value = 900
$ws.Range("H19").Value = 1395

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit) ---
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = 3510
$ws.Range("F20").Value = 7206
$ws.Range("G20").Value = 9996
$ws.Range("H20").Value = 12646

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -172
$ws.Range("E21").Value = -608
$ws.Range("F21").Value = -725
$ws.Range("G21").Value = -1293
$ws.Range("H21").Value = -1160

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم (Net continuing profit) ---
$ws.Range("D22").Value = "-"
$ws.Range("E22").Value = 2902
$ws.Range("F22").Value = 6481
$ws.Range("G22").Value = 8703
$ws.Range("H22").Value = 11486

# --- Row 23: سود (زیان) عملیات متوقف شده (unchanged, still all dashes) ---
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = 2902
$ws.Range("F24").Value = 6481
$ws.Range("G24").Value = 8703
$ws.Range("H24").Value = 11486

# --- Row 25: سود هر سهم پس از کسر مالیات (EPS after tax) ---
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 5437
$ws.Range("E26").Value = 8575
$ws.Range("F26").Value = 4865
$ws.Range("G26").Value = 4169
$ws.Range("H26").Value = 3117

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS on latest capital) ---
$ws.Range("D27").Value = "-"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
